# Pantalla simulador analista - Retanqueo multiple
# Adds new columns (Q:X) of header/data to the "RetanqueoMultiple" sheet,
# mirroring the additional fields captured for the multi-retanqueo flow,
# and widens several columns to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RetanqueoMultiple")
$ws.Activate()

# ---- New header row (row 1) ----
$ws.Range("Q1").Value = "NombreCredito"
$ws.Range("R1").Value = "Mes"
$ws.Range("S1").Value = "fecha"
$ws.Range("T1").Value = "AnnoAfetacion"
$ws.Range("U1").Value = "fechaActual"
$ws.Range("V1").Value = "Banco"
$ws.Range("W1").Value = "Cartera1"
$ws.Range("X1").Value = "Saneamiento2"

# ---- New data row (row 2) ----
# "fechaActual" is written before "NombreCredito" so the shared-string
# table grows in the same order as the source workbook (index 56 =
# "25/10/2021", index 57 = "LUIS CARLOS").
$ws.Range("U2").Value = '"25/10/2021"'
$ws.Range("Q2").Value = '"LUIS CARLOS"'
$ws.Range("R2").Value = '"Octubre"'
$ws.Range("S2").Value = '"14/06/1969"'
$ws.Range("T2").Value = '"2021"'
$ws.Range("V2").Value = '"Remanentes - 60237038927 - REMANENTE"'
$ws.Range("W2").Value = '"0"'
$ws.Range("X2").Value = '"0"'

# ---- Column width adjustments ----
# The runtime quantizes ColumnWidth to 1/6-character steps before writing
# it to the OOXML "width" attribute, so the inputs below are pre-compensated
# (empirically) to land as close as possible on the target widths from the
# source workbook.
$ws.Columns.Item(1).ColumnWidth = 8.752040256892231
$ws.Columns.Item(2).ColumnWidth = 9.752040256892231
$ws.Columns.Item(8).ColumnWidth = 38.41748707706767
$ws.Columns.Item(9).ColumnWidth = 4.417487077067669
$ws.Columns.Item(10).ColumnWidth = 4.75204025689223
$ws.Columns.Item(11).ColumnWidth = 18.752040256892233
$ws.Columns.Item(12).ColumnWidth = 8.752040256892231
$ws.Columns.Item(13).ColumnWidth = 7.75204025689223
$ws.Columns.Item(15).ColumnWidth = 21.584177239974938
$ws.Columns.Item(16).ColumnWidth = 5.91937852443609
$ws.Columns.Item(17).ColumnWidth = 14.25014880952381
$ws.Columns.Item(18).ColumnWidth = 8.919378524436091
$ws.Columns.Item(19).ColumnWidth = 11.584177239974938
$ws.Columns.Item(20).ColumnWidth = 13.417487077067669
$ws.Columns.Item(21).ColumnWidth = 11.584177239974938
$ws.Columns.Item(22).ColumnWidth = 38.75204025689223
$ws.Columns.Item(23).ColumnWidth = 7.584177239974937
$ws.Columns.Item(24).ColumnWidth = 12.752040256892231

# ---- Selection / view state ----
$ws.Range("H10").Select()

Write-Output "edit complete"
